# B_China_POI.conf sheet cleanup:
# Two rows that were flagged with the red "needs review" highlight style
# have been removed from the POI catalog:
#   kind=195000      / catalog=020D   / 外国首都名   (foreign capital name)
#   kind=19020201    / catalog=01ff01 / 普通岛屿     (ordinary island)
# Deleting the rows (rather than just clearing their contents) shifts every
# row below them up, which is what re-numbers the remaining entries from
# A1:C25 down to A1:C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "kind=195000 / catalog=020D / 外国首都名" is row 19 in the original sheet.
$ws.Rows(19).Delete() | Out-Null

# After that delete, everything below shifted up by one, so the row that
# held "kind=19020201 / catalog=01ff01 / 普通岛屿" (originally row 23) is
# now row 22.
$ws.Rows(22).Delete() | Out-Null

# Leave the cursor where the author's last save left it.
$ws.Range("A16").Select() | Out-Null
